# Integrated 2023 and fixed many issues
#
# Applies the semantic edits captured by the authoritative XML diff:
#   1. Rename the worksheet "papers" -> "PLM2023" (the _xlnm._FilterDatabase
#      defined name reference updates automatically with the rename).
#   2. Normalise the keyword-list punctuation for two rows (swap the
#      middle-dot / semicolon separators for plain commas).
#   3. Wrap the text in B2 and grow row 2 to fit it.
#   4. Nudge the column widths to the values recorded after the edit.
#   5. Leave the active selection on D62 (last data cell), matching the
#      saved cursor position in the authored workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the sheet -----------------------------------------------
$ws.Name = "PLM2023"

# --- 2. Fix keyword separators (middle dot / semicolons -> commas) -----
$ws.Range("D56").Value = "Model Based Systems Engineering, RFLP, Requirements- in-Loop, Product Life-cycle Management,Decision Support"
$ws.Range("D62").Value = "Hybrid Production Structures, Flow Production, Job Shop Production, Flexibility, Interdependencies, Transformability, Life Cycle"

# --- 3. Wrap text for B2 and grow the row to fit it ---------------------
$ws.Range("B2").WrapText = $true
$ws.Range("B2").EntireRow.RowHeight = 48

# --- 4. Column width touch-ups ------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 3.33
$ws.Columns.Item(2).ColumnWidth = 18
$ws.Columns.Item(3).ColumnWidth = 124.83
$ws.Columns.Item(4).ColumnWidth = 113.5

# --- 5. Restore the saved selection -------------------------------------
$ws.Range("D62").Select() | Out-Null
